# The workbook's single sheet "Property1" is renamed to "DataNode" as part
# of unifying the DataNode / DataTable / Entity naming scheme.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Restore the last-saved selection recorded for this sheet (cell H33).
$ws.Range("H33").Select()
